$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "[임피던스 제어(2)] 임피던스 제어란?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/100"

$ws.Range("D39").Value = "The 7 Steps of Machine Learning"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/The-7-Steps-of-Machine-Learning-1"

$ws.Range("D51").Value = "[우분투] 우분투(ubuntu) 버전 확인 명령어"
$ws.Range("E51").Value = "https://bskyvision.com/1161"
